# ============================================================================
# Applies the diff: continues the helper ("pomocniczy") write-up file with
# several new paragraphs documenting further work on the project.
# ============================================================================

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The paragraph ending in '...jest drukowany na ekranie.' gains a trailing
#    space run.
# ---------------------------------------------------------------------------
$rEkranie = $d.Paragraphs.Item(18).Range
$rEkranie.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 2) A brand-new paragraph about the git branch mix-up is inserted right
#    after it.
# ---------------------------------------------------------------------------
$rEkranie2 = $d.Paragraphs.Item(18).Range
$rEkranie2.InsertParagraphAfter()
$rPoNapisaniu = $d.Paragraphs.Item(19).Range
$rPoNapisaniu.Text = 'Po napisaniu tej części kodu zauważałam, że jest domyślnie ustawiony branch „main”, więc utworzyłam branch „master”. Branch „main” nie został usunięty, co może być mylące, ponieważ reszta commitów będzie się znajdowała na innej gałęzi. '

# ---------------------------------------------------------------------------
# 3) The existing 'Zeby ' paragraph is continued with the rest of the
#    sentence and the whole paragraph is highlighted yellow.
# ---------------------------------------------------------------------------
$rZeby = $d.Paragraphs.Item(20).Range
$rZeby.InsertAfter('wywołać funkcję w wierszu poleceń bezpośrednio przez polecenie „python”, trzeba mieć ustawione środowisko zmienne.')
$rZebyFull = $d.Paragraphs.Item(20).Range
$rZebyFull.HighlightColorIndex = 7

# ---------------------------------------------------------------------------
# 4) A run of further paragraphs is appended describing the rest of the work.
#    $nextIdx always holds the 1-based Paragraphs index the next new
#    paragraph will land on once InsertParagraphAfter() splits the range.
# ---------------------------------------------------------------------------
$prevIdx = 20

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p0 = $d.Paragraphs.Item($nextIdx).Range
$p0.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p1 = $d.Paragraphs.Item($nextIdx).Range
$p1.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p2 = $d.Paragraphs.Item($nextIdx).Range
$p2.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p3 = $d.Paragraphs.Item($nextIdx).Range
$p3.Text = 'Następnie została dopisana pomocnicza funkcja, która przelicza wartość z radianów na stopnie, minuty, sekundy.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p4 = $d.Paragraphs.Item($nextIdx).Range
$p4.Text = 'Do napisanych definicji została dodana klauzula __name__, na przykładzie tego samego algorytmu.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p5 = $d.Paragraphs.Item($nextIdx).Range
$p5.Text = 'Korzystając z wymienionych wcześniej materiałów, została dodana funkcja XYZ_to_neu, która przyjmuje cztery argumenty: dX, X, Y i Z. Funkcja wykorzystuje te argumenty do obliczenia macierzy transformacji R, która przekształca wektor przesunięcia dX, wyrażony w układzie współrzędnych kartezjańskich z punktem początkowym w (X, Y, Z), na lokalny układ współrzędnych płaszczyzny stycznej zdefiniowany przez wektor normalny w punkcie (X, Y, Z).  I dopisanie do tej funkcji części __name__.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p6 = $d.Paragraphs.Item($nextIdx).Range
$p6.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p7 = $d.Paragraphs.Item($nextIdx).Range
$p7.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p8 = $d.Paragraphs.Item($nextIdx).Range
$p8.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p9 = $d.Paragraphs.Item($nextIdx).Range
$p9.Text = 'Następnym krokiem było napisanie dokumentacji do kodu w pliku README.md, który będzie się wyświetlał na GitHub. Jest to plik tekstowy, który zawiera podstawowe informacje na temat projektu. W naszym przypadku to: nazwa, opis metod transformacji (do czego służy, dane wejściowe i wyjściowe) , sposób instalacji, przykłady używania funkcji i inne ważne informacje, które mogą być przydatne dla użytkowników.'
$uStart = $p9.Start + $p9.Text.IndexOf('sposób instalacji,')
$uEnd = $uStart + 18
$uRange = $d.Range($uStart, $uEnd)
$uRange.Font.Underline = 1
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p10 = $d.Paragraphs.Item($nextIdx).Range
$p10.Text = 'Przy napisaniu wystąpił problem ze śledzeniem plików, ponieważ dodałam do repozytorium plik .png, który wykorzystałam w dokumentacji. Przy rozwiązaniu tego problemu stworzył się w śledzonym folderze clon naszego repozytorium, co powodowało zapisywanie zmian do nowego pliku. Za pomocą  prowadzącego ten problem został rozwiązany, w wyniku czego nie było negatywnych skutków. '
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p11 = $d.Paragraphs.Item($nextIdx).Range
$p11.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p12 = $d.Paragraphs.Item($nextIdx).Range
$p12.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p13 = $d.Paragraphs.Item($nextIdx).Range
$p13.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p14 = $d.Paragraphs.Item($nextIdx).Range
$p14.Text = 'Żeby pisać sprawozdanie z projektu przy pomocy programu LaTeX był stworzony pomocniczy plik, zawierający informacje na temat zadania, czyli cel ćwiczenia, przebieg i tp. Ponieważ nie miałyśmy jeszcze do czynienia z tym programem i chciałyśmy widzieć, jak musi wyglądać ostateczny plik.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p15 = $d.Paragraphs.Item($nextIdx).Range
$p15.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p16 = $d.Paragraphs.Item($nextIdx).Range
$p16.Text = '.'
$prevIdx = $nextIdx

$d.Paragraphs.Item($prevIdx).Range.InsertParagraphAfter()
$nextIdx = $prevIdx + 1
$p17 = $d.Paragraphs.Item($nextIdx).Range
$p17.Text = '.'
$prevIdx = $nextIdx

